# Auto-generated edit script to apply numeric corrections to Sheets
$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 10013.728
$ws.Range("I28").Value = 12837.75
$ws.Range("K28").Value = 12837.75
$ws.Range("M28").Value = -12352.75
$ws.Range("H32").Value = 1076.4445
$ws.Range("I32").Value = 820.7778
$ws.Range("J32").Value = 1332.1111
$ws.Range("K32").Value = 820.7778
$ws.Range("L32").Value = 1332.1111
$ws.Range("M32").Value = -494.7778
$ws.Range("N32").Value = -1984.1111
$ws.Range("H62").Value = 5468.6
$ws.Range("I62").Value = 4655
$ws.Range("K62").Value = 4655
$ws.Range("M62").Value = -4031
$ws.Range("H64").Value = 4999.75
$ws.Range("I64").Value = 4999.75
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 4999.75
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = -4751.75
$ws.Range("N64").ClearContents()
$ws.Range("H65").Value = 5468.6
$ws.Range("I65").Value = 4655
$ws.Range("K65").Value = 23275
$ws.Range("M65").Value = -20155
$ws.Range("H67").Value = 4999.75
$ws.Range("I67").Value = 4999.75
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 4999.75
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = -4141.75
$ws.Range("N67").ClearContents()
$ws.Range("H86").Value = 4799
$ws.Range("I86").Value = 4799
$ws.Range("K86").Value = 4799
$ws.Range("M86").Value = -3676
$ws.Range("H89").Value = 4799
$ws.Range("I89").Value = 4799
$ws.Range("K89").Value = 23995
$ws.Range("M89").Value = -18379
$ws.Range("H98").Value = 1059.7778
$ws.Range("I98").Value = 665
$ws.Range("K98").Value = 665
$ws.Range("M98").Value = 833
$ws.Range("H106").Value = 4998.3335
$ws.Range("I106").Value = 4998.3335
$ws.Range("K106").Value = 4998.3335
$ws.Range("M106").Value = -4367.3335
$ws.Range("H122").Value = 1059.7778
$ws.Range("I122").Value = 665
$ws.Range("K122").Value = 1995
$ws.Range("M122").Value = 455

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5666
$ws.Range("I2").Value = 4499
$ws.Range("K2").Value = 4499
$ws.Range("M2").Value = -4386
$ws.Range("H32").Value = 8665.629000000001
$ws.Range("I32").Value = 3433.2424
$ws.Range("J32").Value = 95000
$ws.Range("K32").Value = 3433.2424
$ws.Range("L32").Value = 95000
$ws.Range("M32").Value = -3146.2424
$ws.Range("N32").Value = -95574
$ws.Range("H45").Value = 1791.0769
$ws.Range("I45").Value = 1732
$ws.Range("K45").Value = 1732
$ws.Range("M45").Value = -1355
$ws.Range("H61").Value = 2363
$ws.Range("I61").Value = 2363
$ws.Range("K61").Value = 2363
$ws.Range("M61").Value = -2151
$ws.Range("H74").Value = 5513
$ws.Range("I74").Value = 4864.1665
$ws.Range("K74").Value = 4864.1665
$ws.Range("M74").Value = -3990.1665
$ws.Range("H77").Value = 5513
$ws.Range("I77").Value = 4864.1665
$ws.Range("K77").Value = 24320.8325
$ws.Range("M77").Value = -19952.8325
$ws.Range("H97").Value = 552.82355
$ws.Range("I97").Value = 552.82355
$ws.Range("K97").Value = 552.82355
$ws.Range("M97").Value = -56.82354999999995
$ws.Range("H102").Value = 2291.5833
$ws.Range("I102").Value = 1659.4
$ws.Range("J102").Value = 5452.5
$ws.Range("K102").Value = 1659.4
$ws.Range("L102").Value = 5452.5
$ws.Range("M102").Value = -37.40000000000009
$ws.Range("N102").Value = -8696.5
$ws.Range("H116").Value = 5666
$ws.Range("I116").Value = 4499
$ws.Range("K116").Value = 4499
$ws.Range("M116").Value = -2205
$ws.Range("H132").Value = 1942.6111
$ws.Range("I132").Value = 1939.2354
$ws.Range("K132").Value = 5817.706200000001
$ws.Range("M132").Value = -3287.706200000001
$ws.Range("H136").Value = 2363
$ws.Range("I136").Value = 2363
$ws.Range("K136").Value = 7089
$ws.Range("M136").Value = -4539

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5666
$ws.Range("I3").Value = 4499
$ws.Range("K3").Value = 4499
$ws.Range("M3").Value = -4385
$ws.Range("H82").Value = 31377.908
$ws.Range("I82").Value = 16289.25
$ws.Range("J82").Value = 40000
$ws.Range("K82").Value = 16289.25
$ws.Range("L82").Value = 40000
$ws.Range("M82").Value = -15906.25
$ws.Range("N82").Value = -40766
$ws.Range("H85").Value = 31377.908
$ws.Range("I85").Value = 16289.25
$ws.Range("J85").Value = 40000
$ws.Range("K85").Value = 16289.25
$ws.Range("L85").Value = 40000
$ws.Range("M85").Value = -14963.25
$ws.Range("N85").Value = -42652
$ws.Range("H94").Value = 2652.4092
$ws.Range("I94").Value = 2617.8
$ws.Range("J94").Value = 2998.5
$ws.Range("K94").Value = 2617.8
$ws.Range("L94").Value = 2998.5
$ws.Range("M94").Value = -2166.8
$ws.Range("N94").Value = -3900.5
$ws.Range("H99").Value = 1173.375
$ws.Range("I99").Value = 1173.375
$ws.Range("K99").Value = 1173.375
$ws.Range("M99").Value = 324.625
$ws.Range("H105").Value = 2709.3333
$ws.Range("I105").Value = 1939
$ws.Range("J105").Value = 4250
$ws.Range("K105").Value = 1939
$ws.Range("L105").Value = 4250
$ws.Range("M105").Value = -192
$ws.Range("N105").Value = -7744

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H43").Value = 9564
$ws.Range("J43").Value = 9564
$ws.Range("L43").Value = 9564
$ws.Range("N43").Value = -9932
$ws.Range("H101").Value = 9564
$ws.Range("J101").Value = 9564
$ws.Range("L101").Value = 9564
$ws.Range("N101").Value = -16054
$ws.Range("H134").Value = 1721.963
$ws.Range("J134").Value = 950
$ws.Range("L134").Value = 2850
$ws.Range("N134").Value = -7920

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 152
$ws.Range("I2").Value = 188.875
$ws.Range("K2").Value = 188.875
$ws.Range("M2").Value = -75.875
$ws.Range("H97").Value = 453
$ws.Range("I97").Value = 270.16666
$ws.Range("K97").Value = 270.16666
$ws.Range("M97").Value = 225.83334
$ws.Range("H101").Value = 22000
$ws.Range("J101").Value = 22000
$ws.Range("L101").Value = 22000
$ws.Range("N101").Value = -28490
$ws.Range("H107").Value = 859.6
$ws.Range("I107").Value = 324.5
$ws.Range("J107").Value = 3000
$ws.Range("K107").Value = 324.5
$ws.Range("L107").Value = 3000
$ws.Range("M107").Value = 1595.5
$ws.Range("N107").Value = -6840
$ws.Range("H113").Value = 2471.0833
$ws.Range("I113").Value = 2471.0833
$ws.Range("K113").Value = 2471.0833
$ws.Range("M113").Value = -301.0832999999998
$ws.Range("H126").Value = 3000
$ws.Range("I126").Value = 2555.5557
$ws.Range("K126").Value = 7666.6671
$ws.Range("M126").Value = -5196.6671

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 314.9375
$ws.Range("I55").Value = 267.7
$ws.Range("K55").Value = 267.7
$ws.Range("M55").Value = -94.69999999999999
$ws.Range("H132").Value = 5068.273
$ws.Range("I132").Value = 2458.5
$ws.Range("K132").Value = 7375.5
$ws.Range("M132").Value = -4845.5

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 2125
$ws.Range("I100").Value = 2125
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 4250
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -3709
$ws.Range("N100").ClearContents()
$ws.Range("H107").Value = 498
$ws.Range("I107").Value = 498
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1494
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 426
$ws.Range("N107").ClearContents()
$ws.Range("H113").Value = 6202.5
$ws.Range("I113").Value = 10425
$ws.Range("K113").Value = 31275
$ws.Range("M113").Value = -29105
$ws.Range("H132").Value = 2203
$ws.Range("I132").Value = 2203
$ws.Range("K132").Value = 6609
$ws.Range("M132").Value = -4079
